$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, centered, thin border) onto the
# new header cells before writing their text, so they match the other
# header cells in row 1 (style index 1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (93 wins, 69 losses, 0 ties) for every data row
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
